$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: columns E/F/G need to widen for the new Orders block ---
$ws.Columns.Item(5).ColumnWidth = 13.7
$ws.Columns.Item(6).ColumnWidth = 16.6
$ws.Columns.Item(7).ColumnWidth = 16.7

# --- Table1: remove Order1, Order2, Order3 columns ---
$lo1 = $ws.ListObjects.Item("Table1")
$lo1.ListColumns.Item("Order3").Delete()
$lo1.ListColumns.Item("Order2").Delete()
$lo1.ListColumns.Item("Order1").Delete()

# --- Table2: add a new drink row (row 20) and resize the table ---
$ws.Cells.Item(20, 1).Value = 7
$ws.Cells.Item(20, 2).Value = "Long Black"
$ws.Cells.Item(20, 3).Value = 3.5

$lo2 = $ws.ListObjects.Item("Table2")
$lo2.Resize($ws.Range("A12:C20"))

# --- New helper / lookup block (rows 23-29, columns A-F) ---
$ws.Range("A23").Value = "Orders"

$ws.Range("A24").Value = "OrderID"
$ws.Range("B24").Value = "PatreonID"
$ws.Range("C24").Value = "DrinkID"

$ws.Range("E23").Value = "OrderPatreons"
$ws.Range("E24").Value = "OrderID"
$ws.Range("F24").Value = "PatreonID"

$orders = @(
    @(1, 1, 1),
    @(2, 1, 2),
    @(3, 1, 5),
    @(4, 2, 2),
    @(5, 3, 3),
    @(6, 3, 5),
    @(7, 4, 7),
    @(8, 4, 4),
    @(9, 5, 4),
    @(10, 5, 6)
)
for ($i = 0; $i -lt $orders.Length; $i++) {
    $r = 25 + $i
    $ws.Cells.Item($r, 1).Value = $orders[$i][0]
    $ws.Cells.Item($r, 2).Value = $orders[$i][1]
    $ws.Cells.Item($r, 3).Value = $orders[$i][2]
}

$patreons = @(1, 2, 3, 4, 5)
for ($i = 0; $i -lt $patreons.Length; $i++) {
    $r = 25 + $i
    $ws.Cells.Item($r, 5).Value = $patreons[$i]
}

# --- New Table3 over the Orders block ---
$lo3 = $ws.ListObjects.Add(1, $ws.Range("A24:C34"), [System.Reflection.Missing]::Value, 1)
$lo3.Name = "Table3"
$lo3.TableStyle = "TableStyleLight10"

# --- Selection / active cell update ---
$ws.Range("F25").Select() | Out-Null
